# Actualizando archivos .xlsx para cargas distribuidas y resortes
#
# Insert three new leading node-id columns (nodo_i, nodo_j, nodo_k) on the
# "carga_distr" (distributed loads) and "kWinkler" (springs) sheets, right
# after the existing "elemento"/"lado" columns and before the traction /
# stiffness component columns, shifting the latter three columns to the
# right.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# carga_distr: elemento, lado, tix, tiy, tjx, tjy, tkx, tky  (A..H)
#  -> elemento, lado, nodo_i, nodo_j, nodo_k, tix, tiy, tjx, tjy, tkx, tky (A..K)
# ---------------------------------------------------------------------
$wsDistr = $wb.Worksheets.Item("carga_distr")

# Drop the unit comments that sat on the columns which are about to move;
# the ones currently on F1/G1/H1 stay put (column insert shifts cell
# CONTENT, not the comment anchors) and will already be correct afterwards.
$wsDistr.Range("C1").Comment.Delete()
$wsDistr.Range("D1").Comment.Delete()
$wsDistr.Range("E1").Comment.Delete()

$wsDistr.Range("C1:E1").EntireColumn.Insert()
$wsDistr.Range("C1").Value = "nodo_i"
$wsDistr.Range("D1").Value = "nodo_j"
$wsDistr.Range("E1").Value = "nodo_k"

$wsDistr.Range("I1").AddComment("N/m^2")

# ---------------------------------------------------------------------
# kWinkler: elemento, lado, kix, kiy, kjx, kjy  (A..F)
#  -> elemento, lado, nodo_i, nodo_j, nodo_k, kix, kiy, kjx, kjy (A..I)
# ---------------------------------------------------------------------
$wsWinkler = $wb.Worksheets.Item("kWinkler")

$wsWinkler.Range("C1").Comment.Delete()
$wsWinkler.Range("D1").Comment.Delete()
$wsWinkler.Range("E1").Comment.Delete()

$wsWinkler.Range("C1:E1").EntireColumn.Insert()
$wsWinkler.Range("C1").Value = "nodo_i"
$wsWinkler.Range("D1").Value = "nodo_j"
$wsWinkler.Range("E1").Value = "nodo_k"

$wsWinkler.Range("G1").AddComment("N/m^2")
$wsWinkler.Range("H1").AddComment("N/m^2")
$wsWinkler.Range("I1").AddComment("N/m^2")

# kWinkler becomes the active sheet/tab.
$wsWinkler.Activate()
$wsWinkler.Range("A1").Select()
